$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 100 into a new row 101 (keeps all of row 100's existing
# values/formatting intact as the appended weekly record), then revise
# rows 99 and 100 in place with the new week's figures.
$ws.Rows.Item(100).Copy()
$ws.Rows.Item(101).Insert()

# --- Row 99: price/date/origin update (existing entry revised) ---
$ws.Range("D99").Value = 45121
$ws.Range("J99").Value = 200
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = 1500
$ws.Range("O99").Value = "Región de Ñuble"
$ws.Range("P99").Value = 1500

# --- Row 100: quality/price update (existing entry revised) ---
$ws.Range("I100").Value = "Primera"
$ws.Range("K100").Value = 1200
$ws.Range("L100").Value = 1200
$ws.Range("M100").Value = 1200
$ws.Range("P100").Value = 1200

# Row 101 already holds the former row 100 data (via the row insert/shift),
# which matches the new weekly record appended to the table.
